# Redeem points 71277620 76.0
#
# 1. Row 17, column A ("phone") was stored as a text value; normalize it
#    back to a genuine number, matching the rest of the column.
# 2. Append a new redemption row (18) for phone 71277620, 76 points,
#    recorded at 2025-08-18T17:10:20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 17: A17 should be numeric, not text ---
$ws.Cells.Item(17, 1).Value = 71277620

# --- Add new row 18 ---
# A18 keeps the phone number as text (matches source data quirk),
# so force a text number format before assigning the numeric-looking string.
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "71277620"

$ws.Cells.Item(18, 2).Value = 76

$ws.Cells.Item(18, 3).Value = "2025-08-18T17:10:20"
